# Applies the 'Updated cryptos list' data refresh (GitHub Actions job).
# Every edited cell is plain text (inlineStr) in the source file, so values
# that look like plain numbers are written with a leading apostrophe to force
# Excel to keep them as text instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.688.13'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '2.474.75'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''319.23'
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').Value = '''92.49'
$ws.Range('E6').Value = '  +1.00%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('D10').Value = '''0.0869'
$ws.Range('E10').Value = '  +8.65%  '
$ws.Range('D11').Value = '''33.21'
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = '2.855.45'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = '''6.89'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').Value = '''15.55'
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('D16').Value = '2.476.05'
$ws.Range('E16').Value = '  +1.40%  '
$ws.Range('D17').Value = '''0.792'
$ws.Range('E17').Value = '  +2.89%  '
$ws.Range('D18').Value = '41.605.64'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').Value = '''6.46'
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').Value = '0.0₃0946'
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').Value = '''70.73'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = '''11.29'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('D23').Value = '''240.12'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('E24').Value = '  +1.77%  '
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '''24.91'
$ws.Range('E27').Value = '  +2.73%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').Value = '''9.72'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('D30').Value = '''36.47'
$ws.Range('E30').Value = '  +4.14%  '
$ws.Range('D31').Value = '''156.91'
$ws.Range('E31').Value = '  +0.75%  '
$ws.Range('D32').Value = '''5.46'
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  +1.07%  '
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('D36').Value = '''17.29'
$ws.Range('E36').Value = '  -0.62%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '''1.85'
$ws.Range('E37').Value = '  +4.12%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = '''0.116'
$ws.Range('E38').Value = '  +1.72%  '
$ws.Range('E39').Value = '  +0.57%  '
$ws.Range('E40').Value = '  +2.64%  '
$ws.Range('E41').Value = '  +2.23%  '
$ws.Range('E42').Value = '  +3.32%  '
$ws.Range('D43').Value = '1.990.85'
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('D44').Value = '''0.0284'
$ws.Range('E44').Value = '  +0.83%  '
$ws.Range('D45').Value = '''18.84'
$ws.Range('E45').Value = '  +1.10%  '
$ws.Range('D46').Value = '''2.98'
$ws.Range('E46').Value = '  +2.92%  '
$ws.Range('D47').Value = '''9.48'
$ws.Range('E47').Value = '  +6.10%  '
$ws.Range('D48').Value = '2.713.22'
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('D49').Value = '''98.15'
$ws.Range('E49').Value = '  +1.79%  '
$ws.Range('D50').Value = '''75.63'
$ws.Range('E50').Value = '  +5.38%  '
$ws.Range('D51').Value = '''67.16'
$ws.Range('E51').Value = '  +1.18%  '
